$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 191, shifting existing rows 191:212 down to 192:213.
$ws.Rows("191").Insert()

# Populate the newly inserted row 191 with the new data record.
$ws.Range("A191").Value = 5
$ws.Range("B191").Value = "Macroferia Regional de Talca"
$ws.Range("C191").Value = "Maule"
$ws.Range("D191").Value = 45013
$ws.Range("E191").Value = 7
$ws.Range("F191").Value = 100112030
$ws.Range("G191").Value = "Poroto granado"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 200
$ws.Range("K191").Value = 30000
$ws.Range("L191").Value = 30000
$ws.Range("M191").Value = 30000
$ws.Range("N191").Value = "$/saco 25 kilos"
$ws.Range("O191").Value = "Región del Maule"
$ws.Range("P191").Value = 1200
$ws.Range("Q191").Value = 25
$ws.Range("R191").Value = "Hortaliza"
